# Upload new version with timestamp
# Fills in the first sale line of the report (row 7 + its row-8 total cell)
# and refreshes the footer (timestamp / page / developer credit) on row 9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

# ---- Row 7: first transaction line --------------------------------------
# A7 = index number
$ws.Range("A7").Value = 1

# C7:G7 = item name -> must be stored as text
$ws.Range("C7:G7").NumberFormat = "@"
$ws.Range("C7").Value = "PANADOL ADVANCE 500 MG 48 TABLETS"

# H7:K7 = current balance (ratio-style text, e.g. "1:2")
$ws.Range("H7:K7").NumberFormat = "@"
$ws.Range("H7").Value = "1:2"

# L7:M7 = order limit - keep its original numeric display format but store
# the value itself as text (matches the source workbook exactly)
$ws.Range("L7").NumberFormat = "@"
$ws.Range("L7").Value = "1"
$ws.Range("L7").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"

# N7:O7 = price -> text
$ws.Range("N7:O7").NumberFormat = "@"
$ws.Range("N7").Value = "92.00"

# P7 = selling price - keep its original numeric display format but store
# the value itself as text
$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "46.0000"
$ws.Range("P7").NumberFormat = "0.00"

# Q7 = number of transactions (ratio-style text, e.g. "0:2")
$ws.Range("Q7").NumberFormat = "@"
$ws.Range("Q7").Value = "0:2"

# ---- Row 8: total for the selling-price column ---------------------------
$ws.Range("P8").Value = 46

# ---- Row 9: footer (generation timestamp / page / credit) ----------------
$ws.Range("A9").Value = "Thursday, 7 August, 2025 9:22 AM"
$ws.Range("G9").Value = "1/1"
$ws.Range("K9").Value = "developed by : Abdelaziz Talaat"
